# Model Result Averages.xlsx - re-worked forecasting / Precision column rescale
# and view-state refresh (active cell / frozen pane scroll position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column N ("Precision") was being stored as a fraction (0-1). Convert it
#    to a percentage (0-100) for every data row (rows 2 through 121).
# ---------------------------------------------------------------------------
$lastRow = 121
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 14)   # column N
    $current = $cell.Value()
    $cell.Value = $current * 100
}

# ---------------------------------------------------------------------------
# 2. Refresh the window view: scroll the frozen pane back up to the top
#    (topLeftCell A2) and move the active selection to S7.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("S7").Select() | Out-Null
